# "Generate Report for Handback"
#
# The localization-status report is regenerated after a handback: the
# "Status" column moves from "Ready for handoff" to
# "Handed back: in sync with en-US", the "Latest Target File" /
# "Latest Handback File" / "Latest Handback DateTime" columns for each
# language sheet are filled in, a hyperlink is added on the new
# "Latest Target File" cell, and a few columns are widened so the new,
# longer values are readable.

$wb = $excel.ActiveWorkbook

$githubUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/501147d9e309f4b6c7e21dee3849dd8e7fd4d96f/e2e/5618eef6-2572-4309-abff-b8645fe5ce31.md"
$mdName = "5618eef6-2572-4309-abff-b8645fe5ce31.md"

# ---------------------------------------------------------------------------
# Overview sheet: widen the per-language status columns (E, F) and flip the
# status text for both languages.
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 29.144371396019366
$overview.Columns.Item(6).ColumnWidth = 29.144371396019366
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 29.144371396019366
$zhcn.Columns.Item(9).ColumnWidth = 39.166666666666664
$zhcn.Columns.Item(10).ColumnWidth = 39.166666666666664

$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("K2").Value = "2016-08-15 16:55:31"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $githubUrl, "", "", $mdName)
$zhcn.Range("I2").Font.Name = "Calibri"
$zhcn.Range("I2").Font.Size = 11
$zhcn.Range("I2").Font.Underline = 2
$zhcn.Range("I2").Font.Color = 15570276

$zhcn.Range("J2").Value = "5618eef6-2572-4309-abff-b8645fe5ce31.eedaa5a4e5c001da6d6e901a393d35f43e221077.zh-cn.xlf"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 29.144371396019366
$dede.Columns.Item(9).ColumnWidth = 39.166666666666664
$dede.Columns.Item(10).ColumnWidth = 39.166666666666664

$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("K2").Value = "2016-08-15 16:55:38"

$dede.Hyperlinks.Add($dede.Range("I2"), $githubUrl, "", "", $mdName)
$dede.Range("I2").Font.Name = "Calibri"
$dede.Range("I2").Font.Size = 11
$dede.Range("I2").Font.Underline = 2
$dede.Range("I2").Font.Color = 15570276

$dede.Range("J2").Value = "5618eef6-2572-4309-abff-b8645fe5ce31.eedaa5a4e5c001da6d6e901a393d35f43e221077.de-de.xlf"
